# Agrego funcionalidad de reporte de clientes
# - Adds 6 new sales rows (rows 19-24) to "Detalle Ventas" for
#   "Torta Argenta (Entera)" sold on 27-02-2025 via WhatsApp (x3) and
#   Instagram (x3), pushing the TOTAL row from 19 -> 25.
# - Updates the TOTAL row's SUBTOTAL formulas to cover the new range.
# - Extends the AutoFilter / _FilterDatabase defined name to A1:J24.
# - Updates the "Medio de Venta" summary counts for WhatsApp/Instagram.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detalle Ventas")

# Insert 6 new rows right above the existing TOTAL row (row 19), which
# pushes the TOTAL row (and its formulas/styles) down to row 25 while
# inheriting formatting from the row above (keeps J column's percent
# style, leaves the rest unstyled) just like the new data rows need.
$ws.Range("A19:A24").EntireRow.Insert()

$newRows = @(
    @("Torta Argenta (Entera)", "27-02-2025", 8, 1, "WhatsApp",   22936.83, 47400, 0, 24463.17, 1.066545377020277),
    @("Torta Argenta (Entera)", "27-02-2025", 8, 1, "WhatsApp",   22936.83, 47400, 0, 24463.17, 1.066545377020277),
    @("Torta Argenta (Entera)", "27-02-2025", 8, 1, "WhatsApp",   22936.83, 47400, 0, 24463.17, 1.066545377020277),
    @("Torta Argenta (Entera)", "27-02-2025", 9, 1, "Instagram",  22936.83, 47400, 0, 24463.17, 1.066545377020277),
    @("Torta Argenta (Entera)", "27-02-2025", 9, 1, "Instagram",  22936.83, 47400, 0, 24463.17, 1.066545377020277),
    @("Torta Argenta (Entera)", "27-02-2025", 9, 1, "Instagram",  22936.83, 47400, 0, 24463.17, 1.066545377020277)
)

$r = 19
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Temporarily blank out the TOTAL row (now row 25) contents so it is not
# contiguous with the data above — otherwise applying AutoFilter below
# would auto-expand to swallow row 25 too (Excel auto-extends the filter
# range to the full contiguous block). Styles on row 25 are untouched.
$ws.Range("A25:J25").ClearContents()

# Re-apply the AutoFilter over the expanded range (toggle off first,
# since the sheet already has one active over the old range).
$ws.AutoFilterMode = $false
$ws.Range("A1:J24").AutoFilter()

# Restore the TOTAL row label + formulas, now pointing at the new range.
$ws.Range("E25").Value = "TOTAL:"
$ws.Range("F25").Formula = "=SUBTOTAL(9,F2:F24)"
$ws.Range("G25").Formula = "=SUBTOTAL(9,G2:G24)"
$ws.Range("H25").Formula = "=SUBTOTAL(9,H2:H24)"
$ws.Range("I25").Formula = "=SUBTOTAL(9,I2:I24)"
$ws.Range("J25").Formula = "=I25/F25"

# Keep the _FilterDatabase defined name in sync with the new AutoFilter range.
$wb.Names.Item("Detalle Ventas!_FilterDatabase").RefersTo = "='Detalle Ventas'!`$A`$1:`$J`$24"

# Update the "Medio de Venta" summary sheet counts.
$ws2 = $wb.Worksheets.Item("Medio de Venta")
$ws2.Range("B3").Value = 5
$ws2.Range("B4").Value = 7
